$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "contains the "
$ws.Range("E2").Value = "yes"
$ws.Range("D2").Select()
